$p = $ppt.ActivePresentation

$oldText = "Generated by AI Automation - Python | 2025-04-24"
$newText = "Generated by AI Automation - Python | 2025-04-25"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $oldHeight = $shp.Height
                $tr.Text = $newText
                $shp.Height = $oldHeight
            }
        }
    }
}
